$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diagnostic Scenarios")
$ws.Activate()

# Insert a new column before the existing "Spawners" column (F) and label it
# "Spawners Current Scenario" - this shifts Spawners/Restoration Potential/
# Rank/Rank Weighted one column to the right (F->G, G->H, H->I, I->J).
$ws.Columns("F").Insert()
$ws.Range("F1").Value = "Spawners Current Scenario"

# Restore the header cell's view/selection state to match the new scroll
# position used while reviewing the added column.
$ws.Range("F15").Select()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1

Write-Output "done"
